$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(4, 2).Value = "cao_tuổi"
$ws.Cells.Item(5, 2).Value = "chăm_sóc"
$ws.Cells.Item(6, 2).Value = "sức_khỏe"
$ws.Cells.Item(7, 2).Value = "xã_hội"
$ws.Cells.Item(8, 2).Value = "tổ_chức"
$ws.Cells.Item(9, 2).Value = "tỉnh"
$ws.Cells.Item(12, 2).Value = "việt_nam"
$ws.Cells.Item(13, 2).Value = "gia_đình"
$ws.Cells.Item(14, 2).Value = "tham_gia"
$ws.Cells.Item(15, 2).Value = "bệnh"
$ws.Cells.Item(17, 2).Value = "huyện"
$ws.Cells.Item(18, 2).Value = "xã"
$ws.Cells.Item(19, 2).Value = "hội_viên"
$ws.Cells.Item(20, 2).Value = "clb"
$ws.Cells.Item(21, 2).Value = "sống"
$ws.Cells.Item(22, 2).Value = "giúp"
$ws.Cells.Item(23, 2).Value = "cụ"
$ws.Cells.Item(24, 2).Value = "công_tác"
$ws.Cells.Item(25, 2).Value = "phát_triển"
$ws.Cells.Item(26, 2).Value = "triệu"
$ws.Cells.Item(27, 2).Value = "vận_động"
$ws.Cells.Item(28, 2).Value = "địa_phương"
$ws.Cells.Item(29, 2).Value = "cơ_sở"
$ws.Cells.Item(32, 2).Value = "phối_hợp"
$ws.Cells.Item(33, 2).Value = "y_tế"
$ws.Cells.Item(34, 2).Value = "kinh_tế"
$ws.Cells.Item(35, 2).Value = "chương_trình"
$ws.Cells.Item(36, 2).Value = "khám"
$ws.Cells.Item(37, 2).Value = "phát_huy"
$ws.Cells.Item(38, 2).Value = "đi"
$ws.Cells.Item(39, 2).Value = "chính_sách"
$ws.Cells.Item(40, 2).Value = "phong_trào"
$ws.Cells.Item(41, 2).Value = "bảo_vệ"
$ws.Cells.Item(42, 2).Value = "dân_số"
$ws.Cells.Item(43, 2).Value = "lao_động"
$ws.Cells.Item(44, 2).Value = "ban"
$ws.Cells.Item(45, 2).Value = "dân"
$ws.Cells.Item(46, 2).Value = "già"
$ws.Cells.Item(47, 2).Value = "văn_hóa"
$ws.Cells.Item(48, 2).Value = "đảng"
$ws.Cells.Item(49, 2).Value = "toàn"
$ws.Cells.Item(50, 2).Value = "cộng_đồng"
$ws.Cells.Item(51, 2).Value = "phường"
$ws.Cells.Item(52, 2).Value = "quà"
$ws.Cells.Item(53, 2).Value = "tích_cực"
$ws.Cells.Item(54, 2).Value = "thành_phố"
$ws.Cells.Item(55, 2).Value = "góp_phần"
$ws.Cells.Item(56, 2).Value = "tinh_thần"
$ws.Cells.Item(57, 2).Value = "triển_khai"
$ws.Cells.Item(58, 2).Value = "mô_hình"
$ws.Cells.Item(59, 2).Value = "hiệu_quả"
$ws.Cells.Item(60, 2).Value = "tặng"
$ws.Cells.Item(61, 2).Value = "hàng"
$ws.Cells.Item(62, 2).Value = "chính_quyền"
$ws.Cells.Item(63, 2).Value = "nghèo"
$ws.Cells.Item(64, 2).Value = "địa_bàn"
$ws.Cells.Item(65, 2).Value = "cán_bộ"
$ws.Cells.Item(66, 2).Value = "con_cháu"
$ws.Cells.Item(67, 2).Value = "tuyên_truyền"
$ws.Cells.Item(68, 2).Value = "cuộc_sống"
$ws.Cells.Item(69, 2).Value = "dinh_dưỡng"
$ws.Cells.Item(70, 2).Value = "ubnd"
$ws.Cells.Item(71, 2).Value = "quỹ"
$ws.Cells.Item(72, 2).Value = "nhà_nước"
$ws.Cells.Item(73, 2).Value = "ngành"
$ws.Cells.Item(74, 2).Value = "đóng_góp"
$ws.Cells.Item(76, 2).Value = "sản_phẩm"
$ws.Cells.Item(77, 2).Value = "quy_định"
$ws.Cells.Item(78, 2).Value = "hoàn_cảnh"
$ws.Cells.Item(79, 2).Value = "môi_trường"
$ws.Cells.Item(80, 2).Value = "trung_tâm"
$ws.Cells.Item(81, 2).Value = "thành_viên"
$ws.Cells.Item(82, 2).Value = "chất_lượng"
$ws.Cells.Item(83, 2).Value = "chủ_tịch"
$ws.Cells.Item(84, 2).Value = "hệ_thống"
$ws.Cells.Item(85, 2).Value = "trung_ương"
$ws.Cells.Item(86, 2).Value = "dịch_vụ"
$ws.Cells.Item(87, 2).Value = "bệnh_viện"
$ws.Cells.Item(88, 2).Value = "nhân_dân"
$ws.Cells.Item(89, 2).Value = "thường_xuyên"
$ws.Cells.Item(90, 2).Value = "điều_trị"
$ws.Cells.Item(91, 2).Value = "trợ_cấp"
$ws.Cells.Item(92, 2).Value = "đường"
$ws.Cells.Item(93, 2).Value = "đời_sống"
$ws.Cells.Item(94, 2).Value = "gương"
$ws.Cells.Item(95, 2).Value = "nhiệm_vụ"
$ws.Cells.Item(96, 2).Value = "thuốc"
$ws.Cells.Item(97, 2).Value = "hằng"
$ws.Cells.Item(98, 2).Value = "nông_thôn"
$ws.Cells.Item(99, 2).Value = "vui"
$ws.Cells.Item(100, 2).Value = "giai_đoạn"
$ws.Cells.Item(101, 2).Value = "hộ"
